$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:E21")
$keyRange = $ws.Range("E2:E21")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2  # xlNo
$ws.Sort.Apply()

$ws.Range("E7").Select() | Out-Null
